# "se termino bien los syllabus"
# Fills in column J (attendance marker "p") for every student row (3-22),
# mirroring the existing column I markers, and updates the sheet's
# active selection to the new last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 22; $row++) {
    $ws.Cells.Item($row, 10).Value = "p"
}

# Leave the selection where the edit naturally ends up (as in the source file)
$ws.Range("J11").Select() | Out-Null
